# Scheduled runner update: refresh market-board derived Leve profit figures
# across the ALC/ARM/BSM/CUL/GSM/LTW/WVR sheets (source data snapshot refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 85867.07000000001
$ws.Range("I132").Value = 99218.25
$ws.Range("K132").Value = 297654.75
$ws.Range("M132").Value = -295124.75
$ws.Range("H133").Value = 85393.71000000001
$ws.Range("J133").Value = 85393.71000000001
$ws.Range("L133").Value = 85393.71000000001
$ws.Range("N133").Value = -95513.71000000001
$ws.Range("H137").Value = 1353663
$ws.Range("I137").Value = 1102.6364
$ws.Range("J137").Value = 2498137.2
$ws.Range("K137").Value = 3307.9092
$ws.Range("L137").Value = 7494411.600000001
$ws.Range("M137").Value = -757.9092000000001
$ws.Range("N137").Value = -7499511.600000001
$ws.Range("H140").Value = 250000
$ws.Range("J140").Value = 250000
$ws.Range("L140").Value = 250000
$ws.Range("N140").Value = -260360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 101.75
$ws.Range("I5").Value = 88.166664
$ws.Range("J5").Value = 142.5
$ws.Range("K5").Value = 88.166664
$ws.Range("L5").Value = 142.5
$ws.Range("M5").Value = 23.833336
$ws.Range("N5").Value = -366.5
$ws.Range("H32").Value = 8134052
$ws.Range("I32").Value = 8134052
$ws.Range("K32").Value = 8134052
$ws.Range("M32").Value = -8133765
$ws.Range("H61").Value = 1668433
$ws.Range("I61").Value = 1755929.5
$ws.Range("K61").Value = 1755929.5
$ws.Range("M61").Value = -1755717.5
$ws.Range("H74").Value = 2762.4333
$ws.Range("I74").Value = 921.5454999999999
$ws.Range("J74").Value = 7824.875
$ws.Range("K74").Value = 921.5454999999999
$ws.Range("L74").Value = 7824.875
$ws.Range("M74").Value = -47.54549999999995
$ws.Range("N74").Value = -9572.875
$ws.Range("H77").Value = 2762.4333
$ws.Range("I77").Value = 921.5454999999999
$ws.Range("J77").Value = 7824.875
$ws.Range("K77").Value = 4607.7275
$ws.Range("L77").Value = 39124.375
$ws.Range("M77").Value = -239.7275
$ws.Range("N77").Value = -47860.375
$ws.Range("H88").Value = 1822.75
$ws.Range("J88").Value = 1396.6666
$ws.Range("L88").Value = 1396.6666
$ws.Range("N88").Value = -2208.6666
$ws.Range("H91").Value = 1822.75
$ws.Range("J91").Value = 1396.6666
$ws.Range("L91").Value = 1396.6666
$ws.Range("N91").Value = -4204.6666
$ws.Range("H97").Value = 1928.5714
$ws.Range("I97").Value = 1750
$ws.Range("K97").Value = 1750
$ws.Range("M97").Value = -1254
$ws.Range("H122").Value = 2771.074
$ws.Range("J122").Value = 3850
$ws.Range("L122").Value = 11550
$ws.Range("N122").Value = -16450
$ws.Range("H136").Value = 1668433
$ws.Range("I136").Value = 1755929.5
$ws.Range("K136").Value = 5267788.5
$ws.Range("M136").Value = -5265238.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 101.75
$ws.Range("I4").Value = 88.166664
$ws.Range("J4").Value = 142.5
$ws.Range("K4").Value = 88.166664
$ws.Range("L4").Value = 142.5
$ws.Range("M4").Value = 26.833336
$ws.Range("N4").Value = -372.5
$ws.Range("H86").Value = 1210.2222
$ws.Range("J86").Value = 1098.25
$ws.Range("L86").Value = 1098.25
$ws.Range("N86").Value = -3344.25
$ws.Range("H89").Value = 1210.2222
$ws.Range("J89").Value = 1098.25
$ws.Range("L89").Value = 5491.25
$ws.Range("N89").Value = -16723.25
$ws.Range("H94").Value = 1200.7273
$ws.Range("I94").Value = 1020.9
$ws.Range("J94").Value = 2999
$ws.Range("K94").Value = 1020.9
$ws.Range("L94").Value = 2999
$ws.Range("M94").Value = -569.9
$ws.Range("N94").Value = -3901
$ws.Range("H116").Value = 65247.668
$ws.Range("J116").Value = 65247.668
$ws.Range("L116").Value = 65247.668
$ws.Range("N116").Value = -74425.66800000001
$ws.Range("H134").Value = 1028568.5
$ws.Range("I134").Value = 1099816.6
$ws.Range("K134").Value = 3299449.8
$ws.Range("M134").Value = -3296914.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 479.5
$ws.Range("I46").Value = 64
$ws.Range("K46").Value = 192
$ws.Range("M46").Value = -101

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3255.721
$ws.Range("I102").Value = 2697.3157
$ws.Range("K102").Value = 2697.3157
$ws.Range("M102").Value = -1075.3157
$ws.Range("H122").Value = 33320.113
$ws.Range("I122").Value = 61407.35
$ws.Range("J122").Value = 6793.278
$ws.Range("K122").Value = 184222.05
$ws.Range("L122").Value = 20379.834
$ws.Range("M122").Value = -181772.05
$ws.Range("N122").Value = -25279.834
$ws.Range("H123").Value = 60001
$ws.Range("J123").Value = 60001
$ws.Range("L123").Value = 60001
$ws.Range("N123").Value = -64901

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4158.278
$ws.Range("I7").Value = 3942.2307
$ws.Range("J7").Value = 4720
$ws.Range("K7").Value = 3942.2307
$ws.Range("L7").Value = 4720
$ws.Range("M7").Value = -3830.2307
$ws.Range("N7").Value = -4944
$ws.Range("H46").Value = 1556
$ws.Range("J46").Value = 2001
$ws.Range("L46").Value = 2001
$ws.Range("N46").Value = -2377
$ws.Range("H68").Value = 4763.5454
$ws.Range("J68").Value = 3624.75
$ws.Range("L68").Value = 3624.75
$ws.Range("N68").Value = -5122.75
$ws.Range("H71").Value = 4763.5454
$ws.Range("J71").Value = 3624.75
$ws.Range("L71").Value = 18123.75
$ws.Range("N71").Value = -25611.75
$ws.Range("H126").Value = 4158.278
$ws.Range("I126").Value = 3942.2307
$ws.Range("J126").Value = 4720
$ws.Range("K126").Value = 11826.6921
$ws.Range("L126").Value = 14160
$ws.Range("M126").Value = -9356.6921
$ws.Range("N126").Value = -19100
$ws.Range("H132").Value = 776849.7
$ws.Range("I132").Value = 1089013.8
$ws.Range("K132").Value = 3267041.4
$ws.Range("M132").Value = -3264511.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1553.7916
$ws.Range("J100").Value = 1954.1428
$ws.Range("L100").Value = 3908.2856
$ws.Range("N100").Value = -4990.2856
$ws.Range("H101").Value = 31845.666
$ws.Range("J101").Value = 31845.666
$ws.Range("L101").Value = 31845.666
$ws.Range("N101").Value = -38335.666
$ws.Range("H105").Value = 80216.164
$ws.Range("J105").Value = 80216.164
$ws.Range("L105").Value = 80216.164
$ws.Range("N105").Value = -87204.164
$ws.Range("H107").Value = 1995.9524
$ws.Range("I107").Value = 906.4545000000001
$ws.Range("J107").Value = 3194.4
$ws.Range("K107").Value = 2719.3635
$ws.Range("L107").Value = 9583.200000000001
$ws.Range("M107").Value = -799.3635000000004
$ws.Range("N107").Value = -13423.2
$ws.Range("H109").Value = 59569
$ws.Range("J109").Value = 59569
$ws.Range("L109").Value = 59569
$ws.Range("N109").Value = -62343
$ws.Range("H126").Value = 706.4
$ws.Range("I126").Value = 706.4
$ws.Range("K126").Value = 2119.2
$ws.Range("M126").Value = 350.8000000000002
$ws.Range("H132").Value = 5300202
$ws.Range("I132").Value = 6492539
$ws.Range("K132").Value = 19477617
$ws.Range("M132").Value = -19475087
$ws.Range("H136").Value = 17091.055
$ws.Range("I136").Value = 10854.4
$ws.Range("K136").Value = 32563.2
$ws.Range("M136").Value = -30013.2

Write-Host "Applied all updates"
